# Scheduled market-data refresh for Belias_Profits workbook
# Updates computed Leve profit figures (current market board averages)
# across all job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5042.857
$ws.Range("I51").Value = 2249.75
$ws.Range("J51").Value = 6160.1
$ws.Range("K51").Value = 2249.75
$ws.Range("L51").Value = 6160.1
$ws.Range("M51").Value = -1765.75
$ws.Range("N51").Value = -7128.1
$ws.Range("H76").Value = 3001.318
$ws.Range("I76").Value = 3000.8125
$ws.Range("K76").Value = 3000.8125
$ws.Range("M76").Value = -2685.8125
$ws.Range("H79").Value = 3001.318
$ws.Range("I79").Value = 3000.8125
$ws.Range("K79").Value = 3000.8125
$ws.Range("M79").Value = -1908.8125
$ws.Range("H126").Value = 24800
$ws.Range("J126").Value = 24800
$ws.Range("L126").Value = 24800
$ws.Range("N126").Value = -34680
$ws.Range("H133").Value = 50500
$ws.Range("J133").Value = 50500
$ws.Range("L133").Value = 50500
$ws.Range("N133").Value = -60620
$ws.Range("H136").Value = 38252.727
$ws.Range("J136").Value = 38252.727
$ws.Range("L136").Value = 38252.727
$ws.Range("N136").Value = -48452.727
$ws.Range("H137").Value = 1183238.1
$ws.Range("I137").Value = 1230.8235
$ws.Range("J137").Value = 4274642
$ws.Range("K137").Value = 3692.4705
$ws.Range("L137").Value = 12823926
$ws.Range("M137").Value = -1142.4705
$ws.Range("N137").Value = -12829026
$ws.Range("H141").Value = 801.55817
$ws.Range("I141").Value = 701.5952
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 2104.7856
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 3075.2144
$ws.Range("N141").Value = -25360

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4870
$ws.Range("I36").Value = 4870
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4870
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4524
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 1335.0834
$ws.Range("I61").Value = 1170.5555
$ws.Range("J61").Value = 1828.6666
$ws.Range("K61").Value = 1170.5555
$ws.Range("L61").Value = 1828.6666
$ws.Range("M61").Value = -958.5554999999999
$ws.Range("N61").Value = -2252.6666
$ws.Range("H88").Value = 3026.6667
$ws.Range("I88").Value = 2980
$ws.Range("J88").Value = 3050
$ws.Range("K88").Value = 2980
$ws.Range("L88").Value = 3050
$ws.Range("M88").Value = -2574
$ws.Range("N88").Value = -3862
$ws.Range("H91").Value = 3026.6667
$ws.Range("I91").Value = 2980
$ws.Range("J91").Value = 3050
$ws.Range("K91").Value = 2980
$ws.Range("L91").Value = 3050
$ws.Range("M91").Value = -1576
$ws.Range("N91").Value = -5858
$ws.Range("H97").Value = 336.08694
$ws.Range("I97").Value = 297.3684
$ws.Range("J97").Value = 520
$ws.Range("K97").Value = 297.3684
$ws.Range("L97").Value = 520
$ws.Range("M97").Value = 198.6316
$ws.Range("N97").Value = -1512
$ws.Range("H132").Value = 2035.6786
$ws.Range("I132").Value = 1908.6111
$ws.Range("K132").Value = 5725.8333
$ws.Range("M132").Value = -3195.8333
$ws.Range("H136").Value = 1335.0834
$ws.Range("I136").Value = 1170.5555
$ws.Range("J136").Value = 1828.6666
$ws.Range("K136").Value = 3511.6665
$ws.Range("L136").Value = 5485.9998
$ws.Range("M136").Value = -961.6664999999998
$ws.Range("N136").Value = -10585.9998

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 39333.332
$ws.Range("J57").Value = 39333.332
$ws.Range("L57").Value = 39333.332
$ws.Range("N57").Value = -40773.332
$ws.Range("H136").Value = 39333.332
$ws.Range("J136").Value = 39333.332
$ws.Range("L136").Value = 39333.332
$ws.Range("N136").Value = -49533.332

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4420.871
$ws.Range("I58").Value = 4921.0386
$ws.Range("K58").Value = 4921.0386
$ws.Range("M58").Value = -4718.0386
$ws.Range("H97").Value = 23500
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 23500
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 23500
$ws.Range("N97").Value = -25482
$ws.Range("M97").ClearContents()
$ws.Range("H122").Value = 1175.3846
$ws.Range("I122").Value = 1074
$ws.Range("J122").Value = 1403.5
$ws.Range("K122").Value = 3222
$ws.Range("L122").Value = 4210.5
$ws.Range("M122").Value = -772
$ws.Range("N122").Value = -9110.5
$ws.Range("H132").Value = 1374555.9
$ws.Range("I132").Value = 2874.7083
$ws.Range("K132").Value = 8624.124899999999
$ws.Range("M132").Value = -6094.124899999999
$ws.Range("H134").Value = 3587.3794
$ws.Range("I134").Value = 3889.348
$ws.Range("J134").Value = 2429.8333
$ws.Range("K134").Value = 11668.044
$ws.Range("L134").Value = 7289.499899999999
$ws.Range("M134").Value = -9133.044
$ws.Range("N134").Value = -12359.4999
$ws.Range("H136").Value = 4420.871
$ws.Range("I136").Value = 4921.0386
$ws.Range("K136").Value = 14763.1158
$ws.Range("M136").Value = -12213.1158

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3615.8333
$ws.Range("J100").Value = 3615.8333
$ws.Range("L100").Value = 10847.4999
$ws.Range("N100").Value = -12469.4999
$ws.Range("H106").Value = 4473.4136
$ws.Range("J106").Value = 4473.4136
$ws.Range("L106").Value = 13420.2408
$ws.Range("N106").Value = -15312.2408

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 23759.8
$ws.Range("J63").Value = 24325
$ws.Range("L63").Value = 24325
$ws.Range("N63").Value = -25697
$ws.Range("H66").Value = 23759.8
$ws.Range("J66").Value = 24325
$ws.Range("L66").Value = 72975
$ws.Range("N66").Value = -79839
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 3033087.8
$ws.Range("I132").Value = 3030.3157
$ws.Range("J132").Value = 7145308.5
$ws.Range("K132").Value = 9090.947100000001
$ws.Range("L132").Value = 21435925.5
$ws.Range("M132").Value = -6560.947100000001
$ws.Range("N132").Value = -21440985.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 770
$ws.Range("I55").Value = 185.25
$ws.Range("J55").Value = 1003.9
$ws.Range("K55").Value = 185.25
$ws.Range("L55").Value = 1003.9
$ws.Range("M55").Value = -12.25
$ws.Range("N55").Value = -1349.9
$ws.Range("H93").Value = 51804.715
$ws.Range("I93").Value = 1337.2
$ws.Range("J93").Value = 177973.5
$ws.Range("K93").Value = 1337.2
$ws.Range("L93").Value = 177973.5
$ws.Range("M93").Value = -89.20000000000005
$ws.Range("N93").Value = -180469.5
$ws.Range("H132").Value = 3164.1572
$ws.Range("I132").Value = 2889.037
$ws.Range("J132").Value = 4092.6875
$ws.Range("K132").Value = 8667.110999999999
$ws.Range("L132").Value = 12278.0625
$ws.Range("M132").Value = -6137.110999999999
$ws.Range("N132").Value = -17338.0625
$ws.Range("H136").Value = 1475.25
$ws.Range("I136").Value = 877.91113
$ws.Range("J136").Value = 3267.2666
$ws.Range("K136").Value = 2633.73339
$ws.Range("L136").Value = 9801.799800000001
$ws.Range("M136").Value = -83.73338999999987
$ws.Range("N136").Value = -14901.7998

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 4500
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = -1797
$ws.Range("N34").Value = -7406
$ws.Range("H132").Value = 2702.1592
$ws.Range("I132").Value = 2827.9688
$ws.Range("J132").Value = 2366.6667
$ws.Range("K132").Value = 8483.9064
$ws.Range("L132").Value = 7100.000100000001
$ws.Range("M132").Value = -5953.9064
$ws.Range("N132").Value = -12160.0001
$ws.Range("H136").Value = 2072.9343
$ws.Range("I136").Value = 1709.5652
$ws.Range("J136").Value = 3187.2666
$ws.Range("K136").Value = 3187.2666
$ws.Range("L136").Value = 9561.799800000001
$ws.Range("M136").Value = -2578.6956
$ws.Range("N136").Value = -14661.7998
